$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 95 and 96 had their match data (columns F:V) swapped - the "Grecia vs
# Liberia" match moved from row 95 to row 96, and the "AD Santos vs Zeledon"
# match moved from row 96 to row 95. Columns A:E (index, country, tournament,
# season, date) stay as they were.
# ---------------------------------------------------------------------------

$ws.Range("F95").Value = "AD Santos"
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = "Zeledon"
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 1.88
$ws.Range("K95").Value = "19/10/2023 18:43"
$ws.Range("L95").Value = 1.97
$ws.Range("M95").Value = "22/10/2023 23:50"
$ws.Range("N95").Value = 3.6
$ws.Range("O95").Value = "19/10/2023 18:43"
$ws.Range("P95").Value = 3.52
$ws.Range("Q95").Value = "22/10/2023 23:50"
$ws.Range("R95").Value = 4.06
$ws.Range("S95").Value = "19/10/2023 18:43"
$ws.Range("T95").Value = 3.95
$ws.Range("U95").Value = "22/10/2023 23:50"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/costa-rica/primera-division/santos-de-guapiles-zeledon/Sb0cnZCg/"

$ws.Range("F96").Value = "Grecia"
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = "Liberia"
$ws.Range("I96").Value = 3
$ws.Range("J96").Value = 2.71
$ws.Range("K96").Value = "19/10/2023 18:43"
$ws.Range("L96").Value = 2.39
$ws.Range("M96").Value = "22/10/2023 23:51"
$ws.Range("N96").Value = 3.38
$ws.Range("O96").Value = "19/10/2023 18:43"
$ws.Range("P96").Value = 3.6
$ws.Range("Q96").Value = "22/10/2023 23:50"
$ws.Range("R96").Value = 2.57
$ws.Range("S96").Value = "19/10/2023 18:43"
$ws.Range("T96").Value = 2.87
$ws.Range("U96").Value = "22/10/2023 23:51"
$ws.Range("V96").Value = "https://www.betexplorer.com/football/costa-rica/primera-division/grecia-liberia/0lhAqXdC/"

# ---------------------------------------------------------------------------
# Append three new match rows (109, 110, 111) at the bottom of the sheet.
# ---------------------------------------------------------------------------

# NOTE: this interpreter only supports POSITIONAL function arguments (named
# "-Param value" args silently bind to nothing), so Set-MatchRow takes its
# values in strict positional order.
function Set-MatchRow {
    param(
        $Sheet, $Row, $Indice, $Data, $Home, $HomeGols, $Away, $AwayGols,
        $HomeOpenOdds, $HomeOpenData, $HomeCloseOdds, $HomeCloseData,
        $DrawOpenOdds, $DrawOpenData, $DrawCloseOdds, $DrawCloseData,
        $AwayOpenOdds, $AwayOpenData, $AwayCloseOdds, $AwayCloseData,
        $Url
    )

    $Sheet.Range("A$Row").Value = $Indice
    $Sheet.Range("B$Row").Value = "costa-rica"
    $Sheet.Range("C$Row").Value = "primera-division"
    $Sheet.Range("D$Row").Value = "2023-2024"
    $Sheet.Range("E$Row").Value = $Data
    $Sheet.Range("F$Row").Value = $Home
    $Sheet.Range("G$Row").Value = $HomeGols
    $Sheet.Range("H$Row").Value = $Away
    $Sheet.Range("I$Row").Value = $AwayGols
    $Sheet.Range("J$Row").Value = $HomeOpenOdds
    $Sheet.Range("K$Row").Value = $HomeOpenData
    $Sheet.Range("L$Row").Value = $HomeCloseOdds
    $Sheet.Range("M$Row").Value = $HomeCloseData
    $Sheet.Range("N$Row").Value = $DrawOpenOdds
    $Sheet.Range("O$Row").Value = $DrawOpenData
    $Sheet.Range("P$Row").Value = $DrawCloseOdds
    $Sheet.Range("Q$Row").Value = $DrawCloseData
    $Sheet.Range("R$Row").Value = $AwayOpenOdds
    $Sheet.Range("S$Row").Value = $AwayOpenData
    $Sheet.Range("T$Row").Value = $AwayCloseOdds
    $Sheet.Range("U$Row").Value = $AwayCloseData
    $Sheet.Range("V$Row").Value = $Url
}

Set-MatchRow $ws 109 108 45238.91666666666 "Zeledon" 0 "Guanacasteca" 0 `
    2.38 "05/11/2023 00:12" 2.9 "08/11/2023 21:53" `
    3.33 "05/11/2023 00:12" 3.3 "08/11/2023 21:53" `
    3.01 "05/11/2023 00:12" 2.53 "08/11/2023 21:53" `
    "https://www.betexplorer.com/football/costa-rica/primera-division/zeledon-guanacasteca/CMIIA9eB/"

Set-MatchRow $ws 110 109 45239.04166666666 "Alajuelense" 2 "Grecia" 0 `
    1.19 "05/11/2023 22:12" 1.25 "09/11/2023 00:56" `
    7.28 "05/11/2023 22:12" 6.33 "09/11/2023 00:57" `
    12.63 "05/11/2023 22:12" 9.99 "09/11/2023 00:57" `
    "https://www.betexplorer.com/football/costa-rica/primera-division/alajuelense-grecia/v9KACVQb/"

Set-MatchRow $ws 111 110 45239.04166666666 "AD Santos" 0 "San Carlos" 2 `
    2.57 "04/11/2023 22:12" 2.51 "09/11/2023 00:59" `
    3.39 "04/11/2023 22:12" 3.37 "09/11/2023 00:56" `
    2.78 "04/11/2023 22:12" 2.87 "09/11/2023 00:59" `
    "https://www.betexplorer.com/football/costa-rica/primera-division/santos-de-guapiles-san-carlos/baC2Eidn/"

# Apply the same formatting as the rest of columns A (bold, centered, thin
# border) and E (datetime number format) to the newly-added rows, by copying
# the format from the row above (row 108) - this reuses the existing style
# entries instead of minting near-duplicate ones.
$xlPasteFormats = -4122
foreach ($r in 109, 110, 111) {
    $ws.Range("A108").Copy()
    $ws.Range("A$r").PasteSpecial($xlPasteFormats)

    $ws.Range("E108").Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteFormats)
}
